$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repulled data: update column F (dSF) values for the listed rows.
$updates = @{
    2  = -1
    3  = -8
    4  = 8
    5  = -4
    6  = -3
    7  = -1
    8  = -1
    9  = -6
    10 = -2
    11 = 5
    12 = -1
    14 = 1
    16 = 1
    18 = -1
    19 = -6
    20 = -1
    21 = -1
    22 = 0
    25 = -7
    26 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
